$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue "D2" '25.869.34'
Set-TextValue "E2" '  -0.44%  '
Set-TextValue "D3" '1.632.96'
Set-TextValue "E3" '  -0.32%  '
Set-TextValue "D4" '1.005'
Set-TextValue "E4" '  +0.32%  '
Set-TextValue "D5" '215.93'
Set-TextValue "E5" '  +0.50%  '
Set-TextValue "D6" '0.5089'
Set-TextValue "D7" '1.006'
Set-TextValue "E7" '  +0.44%  '
Set-TextValue "D8" '0.2590'
Set-TextValue "E8" '  +0.65%  '
Set-TextValue "D9" '0.06411'
Set-TextValue "E9" '  +0.99%  '
Set-TextValue "D10" '19.38'
Set-TextValue "E10" '  -1.93%  '
Set-TextValue "D11" '0.07827'
Set-TextValue "E11" '  +0.74%  '
Set-TextValue "D12" '4.269'
Set-TextValue "E12" '  -0.12%  '
Set-TextValue "D13" '1.629.56'
Set-TextValue "E13" '  -0.45%  '
Set-TextValue "D14" '1.856.60'
Set-TextValue "E14" '  -0.46%  '
Set-TextValue "D15" '0.5597'
Set-TextValue "E15" '  +2.55%  '
Set-TextValue "D16" '63.43'
Set-TextValue "E16" '  -1.35%  '
Set-TextValue "D17" '0.0₅7532'
Set-TextValue "E17" '  -2.61%  '
Set-TextValue "D18" '25.865.40'
Set-TextValue "E19" '  +0.47%  '
Set-TextValue "D20" '193.82'
Set-TextValue "E20" '  -1.43%  '
Set-TextValue "D21" '4.325'
Set-TextValue "E21" '  -2.12%  '
Set-TextValue "D22" '9.847'
Set-TextValue "E22" '  -0.59%  '
Set-TextValue "D23" '6.013'
Set-TextValue "E23" '  -1.09%  '
Set-TextValue "D24" '1.006'
Set-TextValue "E24" '  +0.40%  '
Set-TextValue "D25" '1.836'
Set-TextValue "E25" '  -4.66%  '

# Row 26 and 27 swap: Monero <-> Stellar positions, with updated prices/volumes
Set-TextValue "B26" 'Stellar'
Set-TextValue "C26" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D26" '0.1283'
Set-TextValue "E26" '  +4.36%  '
Set-TextValue "B27" 'Monero'
Set-TextValue "C27" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D27" '140.54'
Set-TextValue "E27" '  -1.48%  '

Set-TextValue "D28" '6.759'
Set-TextValue "E28" '  -1.04%  '
Set-TextValue "D29" '15.44'
Set-TextValue "E29" '  -1.12%  '
Set-TextValue "D30" '1.241'
Set-TextValue "E30" '  +0.19%  '
Set-TextValue "D31" '0.04897'
Set-TextValue "E31" '  +1.07%  '
Set-TextValue "D32" '3.297'
Set-TextValue "E32" '  +0.90%  '
Set-TextValue "D33" '3.190'
Set-TextValue "E33" '  -0.38%  '
Set-TextValue "E34" '  +1.85%  '
Set-TextValue "D35" '2.382'
Set-TextValue "E35" '  +0.34%  '
Set-TextValue "D36" '0.8971'
Set-TextValue "E36" '  -1.68%  '
Set-TextValue "D37" '1.134.78'
Set-TextValue "E37" '  +3.09%  '
Set-TextValue "D38" '2.548'
Set-TextValue "E38" '  -0.80%  '
Set-TextValue "D39" '0.5483'
Set-TextValue "E39" '  -1.01%  '
Set-TextValue "D40" '0.01561'
Set-TextValue "E40" '  -0.42%  '
Set-TextValue "D41" '0.9921'
Set-TextValue "E41" '  -0.89%  '
Set-TextValue "D42" '5.619'
Set-TextValue "E42" '  +1.21%  '
Set-TextValue "D43" '0.7982'
Set-TextValue "E43" '  -0.78%  '
Set-TextValue "E44" '  -1.67%  '
Set-TextValue "D45" '1.779.47'
Set-TextValue "E45" '  +0.05%  '
Set-TextValue "D46" '0.0₈113'
Set-TextValue "E46" '  -6.88%  '
Set-TextValue "D47" '0.4446'
Set-TextValue "E47" '  -1.98%  '
Set-TextValue "E48" '  +0.07%  '
Set-TextValue "D49" '0.05063'
Set-TextValue "E49" '  -2.82%  '
Set-TextValue "D50" '7.591'
Set-TextValue "E50" '  +1.71%  '
Set-TextValue "E51" '  +0.72%  '
